# Auto-generated script to apply market-price / profit-column updates
# captured from the scheduled-runner diff against Sheets/Midgardsormr_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7902.0835
$ws.Range("I19").Value = 7297.6665
$ws.Range("K19").Value = 7297.6665
$ws.Range("M19").Value = -7122.6665
$ws.Range("H31").Value = 6988
$ws.Range("I31").Value = 6988
$ws.Range("K31").Value = 20964
$ws.Range("M31").Value = -20734
$ws.Range("H76").Value = 3876.0625
$ws.Range("I76").Value = 3801.2856
$ws.Range("J76").Value = 4399.5
$ws.Range("K76").Value = 3801.2856
$ws.Range("L76").Value = 4399.5
$ws.Range("M76").Value = -3486.2856
$ws.Range("N76").Value = -5029.5
$ws.Range("H79").Value = 3876.0625
$ws.Range("I79").Value = 3801.2856
$ws.Range("J79").Value = 4399.5
$ws.Range("K79").Value = 3801.2856
$ws.Range("L79").Value = 4399.5
$ws.Range("M79").Value = -2709.2856
$ws.Range("N79").Value = -6583.5
$ws.Range("H97").Value = 2112.25
$ws.Range("J97").Value = 2112.25
$ws.Range("L97").Value = 6336.75
$ws.Range("N97").Value = -7328.75
$ws.Range("H100").Value = 39120.5
$ws.Range("I100").Value = 55448.727
$ws.Range("K100").Value = 55448.727
$ws.Range("M100").Value = -54907.727
$ws.Range("H132").Value = 3882014.2
$ws.Range("I132").Value = 5129081
$ws.Range("K132").Value = 15387243
$ws.Range("M132").Value = -15384713
$ws.Range("H138").Value = 20579.818
$ws.Range("I138").Value = 1667.9744
$ws.Range("J138").Value = 66677.44
$ws.Range("K138").Value = 5003.9232
$ws.Range("L138").Value = 200032.32
$ws.Range("M138").Value = 136.0767999999998
$ws.Range("N138").Value = -210312.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5429.857
$ws.Range("I61").Value = 1456.2727
$ws.Range("J61").Value = 19999.666
$ws.Range("K61").Value = 1456.2727
$ws.Range("L61").Value = 19999.666
$ws.Range("M61").Value = -1244.2727
$ws.Range("N61").Value = -20423.666
$ws.Range("H74").Value = 339824.44
$ws.Range("I74").Value = 353931.75
$ws.Range("K74").Value = 353931.75
$ws.Range("M74").Value = -353057.75
$ws.Range("H77").Value = 339824.44
$ws.Range("I77").Value = 353931.75
$ws.Range("K77").Value = 1769658.75
$ws.Range("M77").Value = -1765290.75
$ws.Range("H88").Value = 6339.364
$ws.Range("I88").Value = 1904.5
$ws.Range("J88").Value = 11661.2
$ws.Range("K88").Value = 1904.5
$ws.Range("L88").Value = 11661.2
$ws.Range("M88").Value = -1498.5
$ws.Range("N88").Value = -12473.2
$ws.Range("H91").Value = 6339.364
$ws.Range("I91").Value = 1904.5
$ws.Range("J91").Value = 11661.2
$ws.Range("K91").Value = 1904.5
$ws.Range("L91").Value = 11661.2
$ws.Range("M91").Value = -500.5
$ws.Range("N91").Value = -14469.2
$ws.Range("H97").Value = 2282.9524
$ws.Range("I97").Value = 1926.1538
$ws.Range("J97").Value = 2862.75
$ws.Range("K97").Value = 1926.1538
$ws.Range("L97").Value = 2862.75
$ws.Range("M97").Value = -1430.1538
$ws.Range("N97").Value = -3854.75
$ws.Range("H102").Value = 5145.909
$ws.Range("I102").Value = 5145.909
$ws.Range("K102").Value = 5145.909
$ws.Range("M102").Value = -3523.909
$ws.Range("H132").Value = 1397.9667
$ws.Range("I132").Value = 1148.2174
$ws.Range("K132").Value = 3444.6522
$ws.Range("M132").Value = -914.6522
$ws.Range("H136").Value = 5429.857
$ws.Range("I136").Value = 1456.2727
$ws.Range("J136").Value = 19999.666
$ws.Range("K136").Value = 4368.8181
$ws.Range("L136").Value = 59998.99800000001
$ws.Range("M136").Value = -1818.8181
$ws.Range("N136").Value = -65098.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1722.7693
$ws.Range("I86").Value = 1789.9
$ws.Range("K86").Value = 1789.9
$ws.Range("M86").Value = -666.9000000000001
$ws.Range("H89").Value = 1722.7693
$ws.Range("I89").Value = 1789.9
$ws.Range("K89").Value = 8949.5
$ws.Range("M89").Value = -3333.5
$ws.Range("H134").Value = 1701.5946
$ws.Range("I134").Value = 1361.5873
$ws.Range("J134").Value = 3648.9092
$ws.Range("K134").Value = 4084.7619
$ws.Range("L134").Value = 10946.7276
$ws.Range("M134").Value = -1549.7619
$ws.Range("N134").Value = -16016.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3127708.8
$ws.Range("I31").Value = 5556607.5
$ws.Range("K31").Value = 5556607.5
$ws.Range("M31").Value = -5556312.5
$ws.Range("H34").Value = 3127708.8
$ws.Range("I34").Value = 5556607.5
$ws.Range("K34").Value = 5556607.5
$ws.Range("M34").Value = -5556405.5
$ws.Range("H107").Value = 909.63635
$ws.Range("J107").Value = 1233.8572
$ws.Range("L107").Value = 1233.8572
$ws.Range("N107").Value = -5073.8572
$ws.Range("H132").Value = 27707.21
$ws.Range("I132").Value = 32120.625
$ws.Range("K132").Value = 96361.875
$ws.Range("M132").Value = -93831.875
$ws.Range("H133").Value = 50099
$ws.Range("I133").Value = 10000
$ws.Range("J133").Value = 63465.332
$ws.Range("K133").Value = 10000
$ws.Range("L133").Value = 63465.332
$ws.Range("M133").Value = -7470
$ws.Range("N133").Value = -68525.33199999999
$ws.Range("H134").Value = 1738.8928
$ws.Range("I134").Value = 1377.9131
$ws.Range("K134").Value = 4133.7393
$ws.Range("M134").Value = -1598.7393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1872.8695
$ws.Range("I107").Value = 3220.3
$ws.Range("J107").Value = 836.38464
$ws.Range("K107").Value = 9660.900000000001
$ws.Range("L107").Value = 2509.15392
$ws.Range("M107").Value = -7740.900000000001
$ws.Range("N107").Value = -6349.15392
$ws.Range("H122").Value = 927.4783
$ws.Range("I122").Value = 529.8570999999999
$ws.Range("J122").Value = 1101.4375
$ws.Range("K122").Value = 4768.7139
$ws.Range("L122").Value = 9912.9375
$ws.Range("M122").Value = -2318.7139
$ws.Range("N122").Value = -14812.9375
$ws.Range("H131").Value = 214138.2
$ws.Range("J131").Value = 2022.7273
$ws.Range("L131").Value = 6068.1819
$ws.Range("N131").Value = -16148.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5679
$ws.Range("I80").Value = 2968.0908
$ws.Range("J80").Value = 8389.909
$ws.Range("K80").Value = 2968.0908
$ws.Range("L80").Value = 8389.909
$ws.Range("M80").Value = -1970.0908
$ws.Range("N80").Value = -10385.909
$ws.Range("H83").Value = 5679
$ws.Range("I83").Value = 2968.0908
$ws.Range("J83").Value = 8389.909
$ws.Range("K83").Value = 14840.454
$ws.Range("L83").Value = 41949.545
$ws.Range("M83").Value = -9848.454
$ws.Range("N83").Value = -51933.545
$ws.Range("H100").Value = 22000
$ws.Range("J100").Value = 22000
$ws.Range("L100").Value = 22000
$ws.Range("N100").Value = -24164
$ws.Range("H102").Value = 5122
$ws.Range("J102").Value = 2166
$ws.Range("L102").Value = 2166
$ws.Range("N102").Value = -5410
$ws.Range("H132").Value = 3003.75
$ws.Range("J132").Value = 1014
$ws.Range("L132").Value = 3042
$ws.Range("N132").Value = -8102

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3469.2942
$ws.Range("I68").Value = 3143.1
$ws.Range("J68").Value = 3935.2856
$ws.Range("K68").Value = 3143.1
$ws.Range("L68").Value = 3935.2856
$ws.Range("M68").Value = -2394.1
$ws.Range("N68").Value = -5433.2856
$ws.Range("H71").Value = 3469.2942
$ws.Range("I71").Value = 3143.1
$ws.Range("J71").Value = 3935.2856
$ws.Range("K71").Value = 15715.5
$ws.Range("L71").Value = 19676.428
$ws.Range("M71").Value = -11971.5
$ws.Range("N71").Value = -27164.428
$ws.Range("H93").Value = 3000.3333
$ws.Range("I93").Value = 3000.3333
$ws.Range("K93").Value = 3000.3333
$ws.Range("M93").Value = -1752.3333
$ws.Range("H122").Value = 3903.138
$ws.Range("I122").Value = 2686.05
$ws.Range("K122").Value = 8058.150000000001
$ws.Range("M122").Value = -5608.150000000001
$ws.Range("H132").Value = 3450.8333
$ws.Range("I132").Value = 3529.762
$ws.Range("K132").Value = 10589.286
$ws.Range("M132").Value = -8059.286
$ws.Range("H136").Value = 2580.682
$ws.Range("I136").Value = 2413.75
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 7241.25
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -4691.25
$ws.Range("N136").Value = -17850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 973.9
$ws.Range("I100").Value = 448.6
$ws.Range("K100").Value = 897.2
$ws.Range("M100").Value = -356.2
$ws.Range("H132").Value = 11674.077
$ws.Range("I132").Value = 12886.811
$ws.Range("K132").Value = 38660.433
$ws.Range("M132").Value = -36130.433
$ws.Range("H141").Value = 91619.875
$ws.Range("J141").Value = 92565.57000000001
$ws.Range("L141").Value = 92565.57000000001
$ws.Range("N141").Value = -102925.57
